# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column inserted
# before the existing "Late" column, pushing "Late" / "heading" /
# "Outstanding" one column to the right (N -> O, O -> P, P -> Q).
# The sheet also becomes the active tab, with cell K13 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet/tab.
$ws.Activate()

# Insert a new blank column before column N ("Late"), shifting the
# remaining columns (Late, heading, Outstanding) one place to the right.
$ws.Columns("N").Insert()

# The newly inserted column inherits the width of the column to its
# left ("In Advance", column M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the selection on K13, matching the saved view state.
$ws.Range("K13").Select() | Out-Null
